$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously had a highlighted PREMIUM cell (yellow fill) because that
# opportunity used to be premium; the refreshed scrape result is not premium,
# so drop the old cell formatting back to the sheet default before writing.
$ws.Range("E2").ClearFormats()

# Keep the opportunity-ID column (A) stored as text for the data rows so
# purely-numeric IDs are not silently auto-converted to numbers.
$ws.Range("A2:A14").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = '1330429'
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1330429'
$ws.Range("C2").Value = 'Machine Learning Intern'
$ws.Range("D2").Value = 'Banur, Punjab, India'
$ws.Range("E2").Value = 'No'
$ws.Range("F2").Value = '1 applicant'
$ws.Range("G2").Value = '3 - 6 Months'
$ws.Range("H2").Value = 'Swami Vivekanand Institute of Engineering & Technology'

# Row 3
$ws.Range("A3").Value = '1330362'
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1330362'
$ws.Range("C3").Value = 'Guest Relations Executive'
$ws.Range("D3").Value = 'Weligama, Sri Lanka'
$ws.Range("E3").Value = 'No'
$ws.Range("F3").Value = '2 applicants'
$ws.Range("G3").Value = '3 - 6 Months'
$ws.Range("H3").Value = 'Steradian Capital Investments'

# Row 4
$ws.Range("A4").Value = '1329520'
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1329520'
$ws.Range("C4").Value = 'Sales Officer'
$ws.Range("D4").Value = 'Dehiwala-Mount Lavinia, Sri Lanka'
$ws.Range("E4").Value = 'No'
$ws.Range("F4").Value = '15 applicants'
$ws.Range("G4").Value = '3 - 6 Months'
$ws.Range("H4").Value = 'Pedro Barn pvt ltd'

# Row 5
$ws.Range("A5").Value = '1329443'
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1329443'
$ws.Range("C5").Value = 'Spa Therapist'
$ws.Range("D5").Value = 'Weligama, Sri Lanka'
$ws.Range("E5").Value = 'No'
$ws.Range("F5").Value = '2 applicants'
$ws.Range("G5").Value = '6 - 18 Months'
$ws.Range("H5").Value = 'Steradian Capital Investments'

# Row 6
$ws.Range("A6").Value = '1328330'
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1328330'
$ws.Range("C6").Value = '[Europe] Accelerate Romania | Digital Marketing & Content Creator'
$ws.Range("D6").Value = 'Brașov, Romania'
$ws.Range("E6").Value = 'No'
$ws.Range("F6").Value = '49 applicants'
$ws.Range("G6").Value = '9 - 12 Weeks'
$ws.Range("H6").Value = 'QHM21 Network'

# Row 7
$ws.Range("A7").Value = '1328041'
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1328041'
$ws.Range("C7").Value = 'IT/ ERP System Developers (Odoo)'
$ws.Range("D7").Value = 'Nugegoda, Sri Lanka'
$ws.Range("E7").Value = 'No'
$ws.Range("F7").Value = '42 applicants'
$ws.Range("G7").Value = '3 - 6 Months'
$ws.Range("H7").Value = 'Epigro Pvt Ltd'

# Row 8
$ws.Range("A8").Value = '1328021'
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1328021'
$ws.Range("C8").Value = 'Guest Relations Officer Intern'
$ws.Range("D8").Value = 'Nugegoda, Sri Lanka'
$ws.Range("E8").Value = 'No'
$ws.Range("F8").Value = '9 applicants'
$ws.Range("G8").Value = '3 - 6 Months'
$ws.Range("H8").Value = 'The Barn By Starbeans in Ella'

# Row 9
$ws.Range("A9").Value = '1327919'
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1327919'
$ws.Range("C9").Value = 'Customer Representative- Intern'
$ws.Range("D9").Value = 'Nugegoda, Sri Lanka'
$ws.Range("E9").Value = 'No'
$ws.Range("F9").Value = '22 applicants'
$ws.Range("G9").Value = '3 - 6 Months'
$ws.Range("H9").Value = 'KAYJAY ELECTRONICS (PVT) LTD'

# Row 10
$ws.Range("A10").Value = '1327813'
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1327813'
$ws.Range("C10").Value = 'Nursery Spanish Practitioner'
$ws.Range("D10").Value = 'Ashby-de-la-Zouch LE65, UK'
$ws.Range("E10").Value = 'No'
$ws.Range("F10").Value = '43 applicants'
$ws.Range("G10").Value = '6 - 18 Months'
$ws.Range("H10").Value = 'Bilingual Day Nursery and Preschool Ltd'

# Row 11
$ws.Range("A11").Value = '1327811'
$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1327811'
$ws.Range("C11").Value = 'Software Engineering Intern'
$ws.Range("D11").Value = 'Colombo, Sri Lanka'
$ws.Range("E11").Value = 'No'
$ws.Range("F11").Value = '116 applicants'
$ws.Range("G11").Value = '3 - 6 Months'
$ws.Range("H11").Value = 'Envision Circle (Pvt) Ltd'

# Row 12
$ws.Range("A12").Value = '1325417'
$ws.Range("B12").Value = 'https://aiesec.org/opportunity/global-talent/1325417'
$ws.Range("C12").Value = 'Junior Software Engineer – AI & Internal Tools (EU ONLY)'
$ws.Range("D12").Value = 'Brussels, Belgium'
$ws.Range("E12").Value = 'No'
$ws.Range("F12").Value = '138 applicants'
$ws.Range("G12").Value = '6 - 18 Months'
$ws.Range("H12").Value = 'Eureka Resource Mining'

# Row 13
$ws.Range("A13").Value = '1325033'
$ws.Range("B13").Value = 'https://aiesec.org/opportunity/global-talent/1325033'
$ws.Range("C13").Value = 'Junior Full-Stack Developer – AI & Web Projects (EU ONLY)'
$ws.Range("D13").Value = 'Brussels, Belgium'
$ws.Range("E13").Value = 'No'
$ws.Range("F13").Value = '155 applicants'
$ws.Range("G13").Value = '6 - 18 Months'
$ws.Range("H13").Value = 'Eureka Resource Mining'

# Row 14
$ws.Range("A14").Value = '1324995'
$ws.Range("B14").Value = 'https://aiesec.org/opportunity/global-talent/1324995'
$ws.Range("C14").Value = '[Partly Remote] Marketing Intern'
$ws.Range("D14").Value = 'Kuala Lumpur, Federal Territory of Kuala Lumpur, Malaysia'
$ws.Range("E14").Value = 'No'
$ws.Range("F14").Value = '56 applicants'
$ws.Range("G14").Value = 'Partly Remote'
$ws.Range("H14").Value = 'Boostorder Sdn. Bhd.'

# Column width adjustments to fit the refreshed content.
# Note: the ColumnWidth property (in characters) reads/writes 5/6 of a
# character lower than the stored OOXML "width" units on this font/theme,
# so subtract 5/6 from the desired stored width to land exactly on target.
$ws.Columns.Item(3).ColumnWidth = 68 - 5/6   # -> stored width 68
$ws.Columns.Item(4).ColumnWidth = 60 - 5/6   # -> stored width 60
$ws.Columns.Item(6).ColumnWidth = 17 - 5/6   # -> stored width 17
$ws.Columns.Item(8).ColumnWidth = 57 - 5/6   # -> stored width 57
